# This workbook was re-saved by a newer version of Excel. Most of the
# resulting diff (new mc:/xr:/x15: namespaces, fileVersion/calcId bumps,
# AlternateContent/revisionPtr/extLst blocks, theme "panose"/typeface
# touch-ups, indexed->auto border colors, x14ac:dyDescent rounding, …) is
# just version metadata that Excel regenerates automatically on every save
# and isn't exposed as an object-model property to script against.
#
# The part of the diff that reflects a genuine user action is the
# view/selection state and the column widths on both sheets, so that's
# what we reproduce here.

$wb = $excel.ActiveWorkbook

# --- "Run Data" sheet: columns A:D narrow slightly, 13.5703125 -> 13.5 chars ---
$wsData = $wb.Worksheets.Item(1)
$wsData.Columns.Item(1).ColumnWidth = 12.67
$wsData.Columns.Item(2).ColumnWidth = 12.67
$wsData.Columns.Item(3).ColumnWidth = 12.67
$wsData.Columns.Item(4).ColumnWidth = 12.67

# --- "Run Info" sheet: col A narrows to 13.5 chars, col B widens (best-fit
# to its longest entry, the long "Comments" string) to ~39.66 chars ---
$wsInfo = $wb.Worksheets.Item(2)
$wsInfo.Columns.Item(1).ColumnWidth = 12.67
$wsInfo.Columns.Item(2).ColumnWidth = 38.75

# --- Selection / active-sheet state ---
# Previously "Run Info" was the active tab with A1 selected. Afterwards
# "Run Info" is left inactive with B12 last selected, and "Run Data"
# becomes the active/selected tab with D5 selected.
$wsInfo.Activate()
$wsInfo.Range("B12").Select()

$wsData.Activate()
$wsData.Range("D5").Select()
